# Modify lobby model; Implement openLobby and closeLobby
#
# The admin screen command is split into an "active" and "inactive" variant
# to reflect the new lobby open/close behaviour:
#   - SHOW_ADMIN_SCREEN  -> SHOW_ACTIVE_ADMIN_SCREEN (row 5), now carrying
#     params ("queue, enqueueKey") and a comment about embedding the
#     enqueueKey in the QR code.
#   - A new SHOW_INACTIVE_ADMIN_SCREEN command is inserted in row 6.
#   - The former row 6 (SHOW_INVALID_INPUT_ERROR_MESSAGE / msg) moves down
#     to row 7, but loses its Params value in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Login / SHOW_ACTIVE_ADMIN_SCREEN / queue, enqueueKey / enqueueKey in QR-Code einbinden
$ws.Range("B5").Value = "SHOW_ACTIVE_ADMIN_SCREEN"
$ws.Range("C5").Value = "queue, enqueueKey "

# Row 6: SHOW_INACTIVE_ADMIN_SCREEN, Params/Comment cleared
$ws.Range("B6").Value = "SHOW_INACTIVE_ADMIN_SCREEN"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""

$ws.Range("D5").Value = "enqueueKey in QR-Code einbinden"

# Row 7: SHOW_INVALID_INPUT_ERROR_MESSAGE, Params/Comment cleared
$ws.Range("B7").Value = "SHOW_INVALID_INPUT_ERROR_MESSAGE"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""

# Widen column C to fit the new "queue, enqueueKey" text (best-fit ~17 chars)
$ws.Columns.Item(3).ColumnWidth = 16.1666666666667

# Update the active cell selection to D8
$ws.Range("D8").Select()
